$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 'Última actualización: 12:01:11'
$ws.Range("A3").Value = 'Total filas: 201'

$ws.Range("A142").Value = '10:13:53'
$ws.Range("B142").Value = '11:20'
$ws.Range("C142").Value = '26_HERNANDEZ'
$ws.Range("D142").Value = 67
$ws.Range("E142").Value = 'LP1912'

$ws.Range("A143").Value = '11:17:39'
$ws.Range("B143").Value = '11:20'
$ws.Range("C143").Value = '225_C ROCA-H SUR'
$ws.Range("D143").Value = 3
$ws.Range("E143").Value = 'LP1912'

$ws.Range("A161").Value = '12:01:11'
$ws.Range("B161").Value = '12:05'
$ws.Range("C161").Value = '23_HERNANDEZ'
$ws.Range("D161").Value = 4
$ws.Range("E161").Value = 'LP1912'

$ws.Range("A162").Value = '10:13:53'
$ws.Range("B162").Value = '12:06'
$ws.Range("C162").Value = '14_ABASTO'
$ws.Range("D162").Value = 113
$ws.Range("E162").Value = 'LP1912'

$ws.Range("A163").Value = '10:13:53'
$ws.Range("B163").Value = '12:06'
$ws.Range("C163").Value = '16_P MOR-SANTA ANA'
$ws.Range("D163").Value = 113
$ws.Range("E163").Value = 'LP1912'

$ws.Range("A164").Value = '10:52:37'
$ws.Range("B164").Value = '12:06'
$ws.Range("C164").Value = '10_OLMOS'
$ws.Range("D164").Value = 74
$ws.Range("E164").Value = 'LP1912'

$ws.Range("A165").Value = '11:46:46'
$ws.Range("B165").Value = '12:07'
$ws.Range("C165").Value = '23_HERNANDEZ'
$ws.Range("D165").Value = 21
$ws.Range("E165").Value = 'LP1912'

$ws.Range("A166").Value = '12:01:11'
$ws.Range("B166").Value = '12:07'
$ws.Range("C166").Value = '16_P MOR-SANTA ANA'
$ws.Range("D166").Value = 6
$ws.Range("E166").Value = 'LP1912'

$ws.Range("A167").Value = '11:17:39'
$ws.Range("B167").Value = '12:13'
$ws.Range("C167").Value = '10_OLMOS'
$ws.Range("D167").Value = 56
$ws.Range("E167").Value = 'LP1912'

$ws.Range("A168").Value = '12:01:11'
$ws.Range("B168").Value = '12:14'
$ws.Range("C168").Value = '17_ROMERO'
$ws.Range("D168").Value = 13
$ws.Range("E168").Value = 'LP1912'

$ws.Range("A169").Value = '10:52:37'
$ws.Range("B169").Value = '12:16'
$ws.Range("C169").Value = '14_ABASTO'
$ws.Range("D169").Value = 84
$ws.Range("E169").Value = 'LP1912'

$ws.Range("A170").Value = '12:01:11'
$ws.Range("B170").Value = '12:17'
$ws.Range("C170").Value = '16_SANTA ANA'
$ws.Range("D170").Value = 16
$ws.Range("E170").Value = 'LP1912'

$ws.Range("A171").Value = '10:52:37'
$ws.Range("B171").Value = '12:20'
$ws.Range("C171").Value = '215A_EL PATO'
$ws.Range("D171").Value = 88
$ws.Range("E171").Value = 'LP1912'

$ws.Range("A172").Value = '11:17:39'
$ws.Range("B172").Value = '12:20'
$ws.Range("C172").Value = '14_ABASTO'
$ws.Range("D172").Value = 63
$ws.Range("E172").Value = 'LP1912'

$ws.Range("A173").Value = '12:01:11'
$ws.Range("B173").Value = '12:21'
$ws.Range("C173").Value = '215A_EL PATO'
$ws.Range("D173").Value = 20
$ws.Range("E173").Value = 'LP1912'

$ws.Range("A174").Value = '12:01:11'
$ws.Range("B174").Value = '12:21'
$ws.Range("C174").Value = '14_ABASTO'
$ws.Range("D174").Value = 20
$ws.Range("E174").Value = 'LP1912'

$ws.Range("A175").Value = '10:52:37'
$ws.Range("B175").Value = '12:21'
$ws.Range("C175").Value = '26_HERNANDEZ'
$ws.Range("D175").Value = 89
$ws.Range("E175").Value = 'LP1912'

$ws.Range("A176").Value = '10:52:37'
$ws.Range("B176").Value = '12:23'
$ws.Range("C176").Value = '17_ROMERO'
$ws.Range("D176").Value = 91
$ws.Range("E176").Value = 'LP1912'

$ws.Range("A177").Value = '12:01:11'
$ws.Range("B177").Value = '12:27'
$ws.Range("C177").Value = '16_SANTA ANA'
$ws.Range("D177").Value = 26
$ws.Range("E177").Value = 'LP1912'

$ws.Range("A178").Value = '11:46:46'
$ws.Range("B178").Value = '12:34'
$ws.Range("C178").Value = '26_HERNANDEZ'
$ws.Range("D178").Value = 48
$ws.Range("E178").Value = 'LP1912'

$ws.Range("A179").Value = '11:17:39'
$ws.Range("B179").Value = '12:34'
$ws.Range("C179").Value = '11_ETCHEVERRY'
$ws.Range("D179").Value = 77
$ws.Range("E179").Value = 'LP1912'

$ws.Range("A180").Value = '12:01:11'
$ws.Range("B180").Value = '12:35'
$ws.Range("C180").Value = '11_ETCHEVERRY'
$ws.Range("D180").Value = 34
$ws.Range("E180").Value = 'LP1912'

$ws.Range("A181").Value = '10:52:37'
$ws.Range("B181").Value = '12:36'
$ws.Range("C181").Value = '27_EL RETIRO'
$ws.Range("D181").Value = 104
$ws.Range("E181").Value = 'LP1912'

$ws.Range("A182").Value = '12:01:11'
$ws.Range("B182").Value = '12:37'
$ws.Range("C182").Value = '27_EL RETIRO'
$ws.Range("D182").Value = 36
$ws.Range("E182").Value = 'LP1912'

$ws.Range("A183").Value = '10:52:37'
$ws.Range("B183").Value = '12:38'
$ws.Range("C183").Value = '17_179 Y 38'
$ws.Range("D183").Value = 106
$ws.Range("E183").Value = 'LP1912'

$ws.Range("A184").Value = '11:46:46'
$ws.Range("B184").Value = '12:41'
$ws.Range("C184").Value = '23_HERNANDEZ'
$ws.Range("D184").Value = 55
$ws.Range("E184").Value = 'LP1912'

$ws.Range("A185").Value = '11:17:39'
$ws.Range("B185").Value = '12:41'
$ws.Range("C185").Value = '10_OLMOS'
$ws.Range("D185").Value = 84
$ws.Range("E185").Value = 'LP1912'

$ws.Range("A186").Value = '12:01:11'
$ws.Range("B186").Value = '12:43'
$ws.Range("C186").Value = '23_HERNANDEZ'
$ws.Range("D186").Value = 42
$ws.Range("E186").Value = 'LP1912'

$ws.Range("A187").Value = '11:17:39'
$ws.Range("B187").Value = '12:48'
$ws.Range("C187").Value = '11_ETCHEVERRY'
$ws.Range("D187").Value = 91
$ws.Range("E187").Value = 'LP1912'

$ws.Range("A188").Value = '12:01:11'
$ws.Range("B188").Value = '12:49'
$ws.Range("C188").Value = '11_ETCHEVERRY'
$ws.Range("D188").Value = 48
$ws.Range("E188").Value = 'LP1912'

$ws.Range("A189").Value = '11:17:39'
$ws.Range("B189").Value = '12:49'
$ws.Range("C189").Value = '17_ROMERO'
$ws.Range("D189").Value = 92
$ws.Range("E189").Value = 'LP1912'

$ws.Range("A190").Value = '10:52:37'
$ws.Range("B190").Value = '12:50'
$ws.Range("C190").Value = '15_ABASTO'
$ws.Range("D190").Value = 118
$ws.Range("E190").Value = 'LP1912'

$ws.Range("A191").Value = '12:01:11'
$ws.Range("B191").Value = '13:02'
$ws.Range("C191").Value = '15_ABASTO'
$ws.Range("D191").Value = 61
$ws.Range("E191").Value = 'LP1912'

$ws.Range("A192").Value = '11:17:39'
$ws.Range("B192").Value = '13:06'
$ws.Range("C192").Value = '16_P MOR-SANTA ANA'
$ws.Range("D192").Value = 109
$ws.Range("E192").Value = 'LP1912'

$ws.Range("A193").Value = '12:01:11'
$ws.Range("B193").Value = '13:07'
$ws.Range("C193").Value = '16_P MOR-SANTA ANA'
$ws.Range("D193").Value = 66
$ws.Range("E193").Value = 'LP1912'

$ws.Range("A194").Value = '11:17:39'
$ws.Range("B194").Value = '13:13'
$ws.Range("C194").Value = '215D_EL PATO'
$ws.Range("D194").Value = 116
$ws.Range("E194").Value = 'LP1912'

$ws.Range("A195").Value = '12:01:11'
$ws.Range("B195").Value = '13:14'
$ws.Range("C195").Value = '215D_EL PATO'
$ws.Range("D195").Value = 73
$ws.Range("E195").Value = 'LP1912'

$ws.Range("A196").Value = '11:46:46'
$ws.Range("B196").Value = '13:17'
$ws.Range("C196").Value = '17_ROMERO'
$ws.Range("D196").Value = 91
$ws.Range("E196").Value = 'LP1912'

$ws.Range("A197").Value = '11:46:46'
$ws.Range("B197").Value = '13:19'
$ws.Range("C197").Value = '10_OLMOS'
$ws.Range("D197").Value = 93
$ws.Range("E197").Value = 'LP1912'

$ws.Range("A198").Value = '12:01:11'
$ws.Range("B198").Value = '13:20'
$ws.Range("C198").Value = '10_OLMOS'
$ws.Range("D198").Value = 79
$ws.Range("E198").Value = 'LP1912'

$ws.Range("A199").Value = '12:01:11'
$ws.Range("B199").Value = '13:21'
$ws.Range("C199").Value = '26_HERNANDEZ'
$ws.Range("D199").Value = 80
$ws.Range("E199").Value = 'LP1912'

$ws.Range("A200").Value = '11:46:46'
$ws.Range("B200").Value = '13:26'
$ws.Range("C200").Value = '14_ABASTO'
$ws.Range("D200").Value = 100
$ws.Range("E200").Value = 'LP1912'

$ws.Range("A201").Value = '11:46:46'
$ws.Range("B201").Value = '13:26'
$ws.Range("C201").Value = '15_ABASTO'
$ws.Range("D201").Value = 100
$ws.Range("E201").Value = 'LP1912'

$ws.Range("A202").Value = '12:01:11'
$ws.Range("B202").Value = '13:27'
$ws.Range("C202").Value = '14_ABASTO'
$ws.Range("D202").Value = 86
$ws.Range("E202").Value = 'LP1912'

$ws.Range("A203").Value = '12:01:11'
$ws.Range("B203").Value = '13:46'
$ws.Range("C203").Value = '17_ROMERO'
$ws.Range("D203").Value = 105
$ws.Range("E203").Value = 'LP1912'

$ws.Range("A204").Value = '12:01:11'
$ws.Range("B204").Value = '13:51'
$ws.Range("C204").Value = '215A_EL PATO'
$ws.Range("D204").Value = 110
$ws.Range("E204").Value = 'LP1912'

$ws.Range("A205").Value = '12:01:11'
$ws.Range("B205").Value = '13:56'
$ws.Range("C205").Value = '225_GOMEZ'
$ws.Range("D205").Value = 115
$ws.Range("E205").Value = 'LP1912'

$ws.Range("A206").Value = '12:01:11'
$ws.Range("B206").Value = '13:57'
$ws.Range("C206").Value = '16_P MOR-167 Y 521'
$ws.Range("D206").Value = 116
$ws.Range("E206").Value = 'LP1912'

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)

$ws.Range("A2").Value = 'Última actualización: 12:01:11'
$ws.Range("A3").Value = 'Total filas: 25'

$ws.Range("A27").Value = '12:01:11'
$ws.Range("B27").Value = '12:21'
$ws.Range("C27").Value = '215A_EL PATO'
$ws.Range("D27").Value = 20
$ws.Range("E27").Value = 'LP1912'

$ws.Range("A28").Value = '11:17:39'
$ws.Range("B28").Value = '13:13'
$ws.Range("C28").Value = '215D_EL PATO'
$ws.Range("D28").Value = 116
$ws.Range("E28").Value = 'LP1912'

$ws.Range("A29").Value = '12:01:11'
$ws.Range("B29").Value = '13:14'
$ws.Range("C29").Value = '215D_EL PATO'
$ws.Range("D29").Value = 73
$ws.Range("E29").Value = 'LP1912'

$ws.Range("A30").Value = '12:01:11'
$ws.Range("B30").Value = '13:51'
$ws.Range("C30").Value = '215A_EL PATO'
$ws.Range("D30").Value = 110
$ws.Range("E30").Value = 'LP1912'

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)

$ws.Range("A2").Value = 'Última actualización: 12:01:11'
$ws.Range("A3").Value = 'Total filas: 33'

$ws.Range("A36").Value = '12:01:11'
$ws.Range("B36").Value = '12:54'
$ws.Range("C36").Value = '215C_LA PLATA'
$ws.Range("D36").Value = 53
$ws.Range("E36").Value = 'L6203'

$ws.Range("A37").Value = '11:46:46'
$ws.Range("B37").Value = '13:30'
$ws.Range("C37").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D37").Value = 104
$ws.Range("E37").Value = 'L6173'

$ws.Range("A38").Value = '12:01:11'
$ws.Range("B38").Value = '13:31'
$ws.Range("C38").Value = '215B_LP-P MOR-1 Y 57'
$ws.Range("D38").Value = 90
$ws.Range("E38").Value = 'L6173'
